# Apply cryptos list update (commit: Updated cryptos list on Sun May 26 12:56:47 UTC 2024 with GitHub Actions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Cells.Item(2, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "68.987.34"
$ws.Cells.Item(2, 5).Value = "  -0.18%  "

$dCell = $ws.Cells.Item(3, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.804.68"
$ws.Cells.Item(3, 5).Value = "  +1.93%  "

$ws.Cells.Item(4, 5).Value = "  +0.02%  "

$dCell = $ws.Cells.Item(5, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "600.82"
$ws.Cells.Item(5, 5).Value = "  -0.19%  "

$dCell = $ws.Cells.Item(6, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "163.44"
$ws.Cells.Item(6, 5).Value = "  -2.62%  "

$dCell = $ws.Cells.Item(7, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.803.33"
$ws.Cells.Item(7, 5).Value = "  +1.95%  "

$ws.Cells.Item(8, 5).Value = "  +0.02%  "

$ws.Cells.Item(9, 5).Value = "  -0.26%  "

$dCell = $ws.Cells.Item(10, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.169"
$ws.Cells.Item(10, 5).Value = "  +1.98%  "

$dCell = $ws.Cells.Item(11, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "6.30"
$ws.Cells.Item(11, 5).Value = "  -1.40%  "

$ws.Cells.Item(12, 5).Value = "  -0.16%  "

$dCell = $ws.Cells.Item(13, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "37.15"
$ws.Cells.Item(13, 5).Value = "  -2.35%  "

$dCell = $ws.Cells.Item(14, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0000245"
$ws.Cells.Item(14, 5).Value = "  -0.42%  "

$dCell = $ws.Cells.Item(15, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "4.443.62"
$ws.Cells.Item(15, 5).Value = "  +2.02%  "

$dCell = $ws.Cells.Item(16, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.799.33"
$ws.Cells.Item(16, 5).Value = "  +1.71%  "

$dCell = $ws.Cells.Item(17, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "69.151.86"
$ws.Cells.Item(17, 5).Value = "  +0.02%  "

$ws.Cells.Item(18, 5).Value = "  +2.43%  "

$ws.Cells.Item(20, 5).Value = "  +5.01%  "

$dCell = $ws.Cells.Item(21, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "17.23"
$ws.Cells.Item(21, 5).Value = "  +1.32%  "

$dCell = $ws.Cells.Item(22, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "485.69"
$ws.Cells.Item(22, 5).Value = "  -1.48%  "

$dCell = $ws.Cells.Item(23, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.720"
$ws.Cells.Item(23, 5).Value = "  -0.65%  "

$dCell = $ws.Cells.Item(24, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.0000161"
$ws.Cells.Item(24, 5).Value = "  +6.47%  "

$dCell = $ws.Cells.Item(25, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "84.70"
$ws.Cells.Item(25, 5).Value = "  +0.02%  "

$ws.Cells.Item(26, 5).Value = "  -2.63%  "

$dCell = $ws.Cells.Item(27, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "12.22"
$ws.Cells.Item(27, 5).Value = "  -0.60%  "

$ws.Cells.Item(28, 5).Value = "  -1.06%  "

$ws.Cells.Item(29, 5).Value = "  -0.08%  "

$ws.Cells.Item(30, 5).Value = "  -1.13%  "

$ws.Cells.Item(31, 5).Value = "  -0.51%  "

$dCell = $ws.Cells.Item(32, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.38"
$ws.Cells.Item(32, 5).Value = "  -4.64%  "

$dCell = $ws.Cells.Item(33, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.955.45"
$ws.Cells.Item(33, 5).Value = "  +1.96%  "

$dCell = $ws.Cells.Item(34, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "31.67"
$ws.Cells.Item(34, 5).Value = "  +0.38%  "

$dCell = $ws.Cells.Item(35, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.750.09"
$ws.Cells.Item(35, 5).Value = "  +2.27%  "

$ws.Cells.Item(36, 5).Value = "  -1.15%  "

$ws.Cells.Item(37, 5).Value = "  +1.66%  "

$dCell = $ws.Cells.Item(38, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.141"
$ws.Cells.Item(38, 5).Value = "  +4.99%  "

$dCell = $ws.Cells.Item(39, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "5.88"
$ws.Cells.Item(39, 5).Value = "  +0.39%  "

$dCell = $ws.Cells.Item(40, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.00"
$ws.Cells.Item(40, 5).Value = "  +0.01%  "

$dCell = $ws.Cells.Item(41, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "0.318"
$ws.Cells.Item(41, 5).Value = "  -1.56%  "

$dCell = $ws.Cells.Item(42, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "3.02"
$ws.Cells.Item(42, 5).Value = "  +1.06%  "

$dCell = $ws.Cells.Item(43, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "429.25"
$ws.Cells.Item(43, 5).Value = "  -0.57%  "

$dCell = $ws.Cells.Item(44, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "48.62"
$ws.Cells.Item(44, 5).Value = "  -0.08%  "

$dCell = $ws.Cells.Item(45, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "1.98"
$ws.Cells.Item(45, 5).Value = "  +0.15%  "

$dCell = $ws.Cells.Item(47, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "8.37"
$ws.Cells.Item(47, 5).Value = "  -1.13%  "

$dCell = $ws.Cells.Item(48, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "2.826.48"
$ws.Cells.Item(48, 5).Value = "  +1.77%  "

$dCell = $ws.Cells.Item(49, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "141.88"
$ws.Cells.Item(49, 5).Value = "  +0.91%  "

$dCell = $ws.Cells.Item(50, 4)
$dCell.NumberFormat = "@"
$dCell.Value = "39.32"
$ws.Cells.Item(50, 5).Value = "  -2.21%  "

$ws.Cells.Item(51, 5).Value = "  -0.07%  "
